$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.754.30"
$ws.Range("E2").Value = "  -0.88%  "
$ws.Range("D3").Value = "1.624.65"
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").Value = "'214.60"
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("D6").Value = "'0.5067"
$ws.Range("E6").Value = "  -1.13%  "
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("D8").Value = "'0.2552"
$ws.Range("E8").Value = "  -1.20%  "
$ws.Range("D9").Value = "'0.06352"
$ws.Range("E9").Value = "  -0.22%  "
$ws.Range("D10").Value = "'19.30"
$ws.Range("E10").Value = "  -2.48%  "
$ws.Range("D11").Value = "'0.07770"
$ws.Range("E11").Value = "  -0.20%  "
$ws.Range("D12").Value = "'4.242"
$ws.Range("E12").Value = "  -1.06%  "
$ws.Range("D13").Value = "1.631.59"
$ws.Range("E13").Value = "  -1.15%  "
$ws.Range("D14").Value = "1.848.79"
$ws.Range("E14").Value = "  -1.13%  "
$ws.Range("D15").Value = "'0.5533"
$ws.Range("E15").Value = "  +1.20%  "
$ws.Range("D16").Value = "'63.61"
$ws.Range("E16").Value = "  -1.47%  "
$ws.Range("D17").Value = "0.0₅7493"
$ws.Range("E17").Value = "  -3.15%  "
$ws.Range("D18").Value = "25.776.57"
$ws.Range("E18").Value = "  -0.95%  "
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("D20").Value = "'4.398"
$ws.Range("E20").Value = "  -0.89%  "
$ws.Range("D21").Value = "'193.52"
$ws.Range("E21").Value = "  -2.12%  "
$ws.Range("D22").Value = "'9.753"
$ws.Range("D23").Value = "'5.967"
$ws.Range("E23").Value = "  -1.91%  "
$ws.Range("D24").Value = "'1.001"
$ws.Range("E24").Value = "  -0.40%  "
$ws.Range("D25").Value = "'1.865"
$ws.Range("E25").Value = "  -1.58%  "
$ws.Range("D26").Value = "'140.71"
$ws.Range("E26").Value = "  -0.82%  "
$ws.Range("D27").Value = "'0.1236"
$ws.Range("E27").Value = "  +0.68%  "
$ws.Range("D28").Value = "'6.725"
$ws.Range("E28").Value = "  -1.99%  "
$ws.Range("D29").Value = "'15.45"
$ws.Range("E29").Value = "  -1.21%  "
$ws.Range("D30").Value = "'1.233"
$ws.Range("E30").Value = "  -0.44%  "
$ws.Range("D31").Value = "'0.04867"
$ws.Range("E31").Value = "  -0.33%  "
$ws.Range("D32").Value = "'3.300"
$ws.Range("E32").Value = "  +0.61%  "
$ws.Range("D33").Value = "'3.174"
$ws.Range("E33").Value = "  -1.05%  "
$ws.Range("D34").Value = "'1.543"
$ws.Range("E34").Value = "  +0.41%  "
$ws.Range("E35").Value = "  -0.62%  "
$ws.Range("D36").Value = "'0.8922"
$ws.Range("E36").Value = "  -2.30%  "
$ws.Range("D37").Value = "1.132.62"
$ws.Range("E37").Value = "  +1.62%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "'0.5499"
$ws.Range("E38").Value = "  -0.48%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'2.536"
$ws.Range("E39").Value = "  -2.11%  "
$ws.Range("D40").Value = "'0.01556"
$ws.Range("E40").Value = "  -0.71%  "
$ws.Range("D41").Value = "'1.001"
$ws.Range("E41").Value = "  -0.30%  "
$ws.Range("D42").Value = "'5.571"
$ws.Range("E42").Value = "  +0.86%  "
$ws.Range("D43").Value = "'0.7926"
$ws.Range("E43").Value = "  -2.21%  "
$ws.Range("D44").Value = "'97.21"
$ws.Range("E44").Value = "  -2.15%  "
$ws.Range("D45").Value = "1.771.65"
$ws.Range("E45").Value = "  -0.43%  "
$ws.Range("D46").Value = "0.0₈114"
$ws.Range("E46").Value = "  -6.45%  "
$ws.Range("D47").Value = "'0.4415"
$ws.Range("E47").Value = "  -2.70%  "
$ws.Range("D48").Value = "'54.64"
$ws.Range("E48").Value = "  -0.83%  "
$ws.Range("D49").Value = "'0.05128"
$ws.Range("E49").Value = "  -2.92%  "
$ws.Range("D50").Value = "'7.588"
$ws.Range("E50").Value = "  +3.30%  "
$ws.Range("D51").Value = "'0.9986"
$ws.Range("E51").Value = "  -0.75%  "
